# Update RNB.Inhab.Year (column C) values for selected countries
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    8 = 53250
    9 = 49350
    12 = 45870
    13 = 1200
    19 = 5700
    21 = 9080
    25 = 44950
    26 = 84450
    37 = 20230
    38 = 47110
    40 = 60170
    44 = 29300
    45 = 21130
    47 = 48300
    48 = 41090
    50 = 41790
    51 = 4450
    56 = 6650
    57 = 19760
    60 = 2320
    63 = 14760
    66 = 61210
    68 = 68120
    69 = 40930
    70 = 33740
    74 = 8080
    78 = 30620
    86 = 17360
    87 = 70840
    88 = 16500
    90 = 3900
    105 = 2020
    107 = 80640
    109 = 41020
    117 = 21980
    130 = 24500
    131 = 55540
    144 = 63170
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

